# ---------------------------------------------------------------------------
# feat: add 2022-Q4 data
#
# 1) "总计" summary sheet gets a new first data row for 2022-Q4 (持有数量(只)=2,
#    持有市值(亿元)=0.04); every other row shifts down by one.
# 2) A brand-new worksheet named "2022-Q4" is inserted right after "总计"
#    (i.e. it becomes the 2nd tab) holding the two funds reported for that
#    quarter. All the other quarter sheets keep their own data untouched and
#    simply slide one tab to the right.
# ---------------------------------------------------------------------------

function Set-TextCell {
    # Forces the value into the cell as TEXT, even when it looks like a
    # number (fund codes such as "009917" must keep their leading zeros,
    # and ratio-looking values like "24.08" must stay strings, matching the
    # source data which stores them as inline strings, not numbers).
    param($Range, [string]$Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: "总计" (summary) sheet - insert the 2022-Q4 row at the top of the
# data (row 2), pushing the existing four rows down and renumbering the
# index column (A).
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Read the existing 4 data rows (row 2..5) before we overwrite anything.
$existingRows = New-Object System.Collections.ArrayList
for ($r = 2; $r -le 5; $r++) {
    $row = @{
        B = $summary.Cells.Item($r, 2).Value()
        C = $summary.Cells.Item($r, 3).Value()
        D = $summary.Cells.Item($r, 4).Value()
    }
    [void]$existingRows.Add($row)
}

# New full data set: 2022-Q4 first, then the previous 4 quarters unchanged.
$newData = New-Object System.Collections.ArrayList
[void]$newData.Add(@{ B = "2022-Q4"; C = 2; D = 0.04 })
foreach ($row in $existingRows) {
    [void]$newData.Add($row)
}

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $i + 2
    $summary.Cells.Item($r, 1).Value = $i
    Set-TextCell $summary.Cells.Item($r, 2) ([string]$newData[$i].B)
    $summary.Cells.Item($r, 3).Value = $newData[$i].C
    $summary.Cells.Item($r, 4).Value = $newData[$i].D
}

# ---------------------------------------------------------------------------
# Step 2: brand-new "2022-Q4" fund sheet, inserted right after "总计".
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $summary)
$q4.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Count; $c++) {
    Set-TextCell $q4.Cells.Item(1, $c + 2) $headers[$c]
}

$funds = @(
    @{ code = "009917"; name = "格林泓利增强债券C"; scale = "2.82"; pos = "24.08"; ratio = "1.17"; mv = "0.0330"; rank = 10 },
    @{ code = "009916"; name = "格林泓利增强债券A"; scale = "0.97"; pos = "24.08"; ratio = "1.17"; mv = "0.0113"; rank = 10 }
)

for ($i = 0; $i -lt $funds.Count; $i++) {
    $r = $i + 2
    $f = $funds[$i]
    $q4.Cells.Item($r, 1).Value = $i
    Set-TextCell $q4.Cells.Item($r, 2) $f.code
    Set-TextCell $q4.Cells.Item($r, 3) $f.name
    Set-TextCell $q4.Cells.Item($r, 4) $f.scale
    Set-TextCell $q4.Cells.Item($r, 5) $f.pos
    Set-TextCell $q4.Cells.Item($r, 6) $f.ratio
    Set-TextCell $q4.Cells.Item($r, 7) $f.mv
    $q4.Cells.Item($r, 8).Value = $f.rank
}

# Restore the originally-selected tab ("2021-Q3", the last sheet) so the
# workbook doesn't open on the newly inserted sheet.
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
